$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "LP1912": new scrape at 01:17:09 — existing row 7 (14_ABASTO)
# gets re-stamped with the new scrape time and a recomputed Minutos,
# and a brand-new row 8 (215_ALUAR / 101) is appended.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:17:09"
$ws1.Range("A3").Value = "Total filas: 3"

$ws1.Range("A7").Value = "01:17:09"
$ws1.Range("D7").Value = 41

$ws1.Range("A8").Value = "01:17:09"
$ws1.Range("B8").Value = "02:58"
$ws1.Range("C8").Value = "215_ALUAR"
$ws1.Range("D8").Value = 101
$ws1.Range("E8").Value = "LP1912"

# ------------------------------------------------------------------
# Sheet "LP1912-215": new scrape at 01:17:09 — appended as row 7.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:17:09"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A7").Value = "01:17:09"
$ws2.Range("B7").Value = "02:58"
$ws2.Range("C7").Value = "215_ALUAR"
$ws2.Range("D7").Value = 101
$ws2.Range("E7").Value = "LP1912"

# ------------------------------------------------------------------
# Sheet "6203-6173": only the "Última actualización" stamp changes.
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:17:09"
